$wb = $excel.ActiveWorkbook

# Excel's serialized column width = ColumnWidth (chars) + 0.8333333333333334 (5/6)
# so to land on an XML width of exactly 40 we need ColumnWidth = 39.166666666666664
$colWidth40 = 39.166666666666664

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e104eae20b290c0ec0cbba8717863275dd923c8c/e2e/ec5acd4f-9c84-40ac-9637-321219b06b1e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52850ea228043ed0720476e588c8840674752fe1/e2e/ec5acd4f-9c84-40ac-9637-321219b06b1e.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52850ea228043ed0720476e588c8840674752fe1/e2e/ec5acd4f-9c84-40ac-9637-321219b06b1e.md"
$displayMd = "ec5acd4f-9c84-40ac-9637-321219b06b1e.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Columns.Item(9).ColumnWidth = $colWidth40
$ws.Columns.Item(10).ColumnWidth = $colWidth40
$ws.Columns.Item(16).ColumnWidth = $colWidth40

$ws.Range("I6").Value = $displayMd
$ws.Range("J6").Value = "ec5acd4f-9c84-40ac-9637-321219b06b1e.8ffc2b31f29f4fc29806e7d53064d2724cca7430.zh-cn.xlf"
$ws.Range("K6").Value = "2016-11-09 00:11:09"
$ws.Range("P6").Value = $errorDetail

$ws.Hyperlinks.Add($ws.Range("I6"), $latestUrl, $null, $null, $displayMd)
$ws.Range("I6").Style = "Hyperlink"

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Columns.Item(9).ColumnWidth = $colWidth40
$ws2.Columns.Item(10).ColumnWidth = $colWidth40
$ws2.Columns.Item(16).ColumnWidth = $colWidth40

$ws2.Range("I6").Value = $displayMd
$ws2.Range("J6").Value = "ec5acd4f-9c84-40ac-9637-321219b06b1e.8ffc2b31f29f4fc29806e7d53064d2724cca7430.de-de.xlf"
$ws2.Range("K6").Value = "2016-11-09 00:11:27"
$ws2.Range("P6").Value = $errorDetail

$ws2.Hyperlinks.Add($ws2.Range("I6"), $latestUrl, $null, $null, $displayMd)
$ws2.Range("I6").Style = "Hyperlink"
